# Scheduled-runner update: refresh cached Universalis market-price
# snapshots (currentAveragePrice*) and the dependent Leve profit
# columns (K:N) on every class sheet of the Leviathan Profits workbook.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 235.18182
$ws.Range("I9").Value = 329.8
$ws.Range("K9").Value = 329.8
$ws.Range("M9").Value = -160.8
# Row 11
$ws.Range("H11").Value = 570.4545000000001
$ws.Range("I11").Value = 570.4545000000001
$ws.Range("K11").Value = 570.4545000000001
$ws.Range("M11").Value = -430.4545000000001
# Row 127
$ws.Range("H127").Value = 1928.8334
$ws.Range("J127").Value = 8000
$ws.Range("L127").Value = 24000
$ws.Range("N127").Value = -33920

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 23164.463
$ws.Range("I32").Value = 4231.936
$ws.Range("J32").Value = 150282.86
$ws.Range("K32").Value = 4231.936
$ws.Range("L32").Value = 150282.86
$ws.Range("M32").Value = -3944.936
$ws.Range("N32").Value = -150856.86
# Row 45
$ws.Range("H45").Value = 6992.4585
$ws.Range("I45").Value = 11282.182
$ws.Range("J45").Value = 3362.6924
$ws.Range("K45").Value = 11282.182
$ws.Range("L45").Value = 3362.6924
$ws.Range("M45").Value = -10905.182
$ws.Range("N45").Value = -4116.6924
# Row 46
$ws.Range("H46").Value = 6801
$ws.Range("I46").Value = 6088
$ws.Range("J46").Value = 7038.6665
$ws.Range("K46").Value = 6088
$ws.Range("L46").Value = 7038.6665
$ws.Range("M46").Value = -5769
$ws.Range("N46").Value = -7676.6665
# Row 61
$ws.Range("H61").Value = 2267.5
$ws.Range("I61").Value = 1932.6129
$ws.Range("K61").Value = 1932.6129
$ws.Range("M61").Value = -1720.6129
# Row 107
$ws.Range("H107").Value = 39999
$ws.Range("J107").Value = 39999
$ws.Range("L107").Value = 39999
$ws.Range("N107").Value = -47679
# Row 136
$ws.Range("H136").Value = 2267.5
$ws.Range("I136").Value = 1932.6129
$ws.Range("K136").Value = 5797.8387
$ws.Range("M136").Value = -3247.8387

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 2030.8235
$ws.Range("I99").Value = 1944
$ws.Range("J99").Value = 2239.2
$ws.Range("K99").Value = 1944
$ws.Range("L99").Value = 2239.2
$ws.Range("M99").Value = -446
$ws.Range("N99").Value = -5235.2
# Row 105
$ws.Range("H105").Value = 1342.6111
$ws.Range("I105").Value = 1297.4615
$ws.Range("K105").Value = 1297.4615
$ws.Range("M105").Value = 449.5385000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1713.409
$ws.Range("I16").Value = 1599.7368
$ws.Range("J16").Value = 2433.3333
$ws.Range("K16").Value = 1599.7368
$ws.Range("L16").Value = 2433.3333
$ws.Range("M16").Value = -1312.7368
$ws.Range("N16").Value = -3007.3333
# Row 93
$ws.Range("H93").Value = 24648.5
$ws.Range("I93").Value = 19998.334
$ws.Range("K93").Value = 19998.334
$ws.Range("M93").Value = -18126.334
# Row 105
$ws.Range("H105").Value = 2573.875
$ws.Range("I105").Value = 2598.8572
$ws.Range("J105").Value = 2399
$ws.Range("K105").Value = 2598.8572
$ws.Range("L105").Value = 2399
$ws.Range("M105").Value = -851.8571999999999
$ws.Range("N105").Value = -5893
# Row 113
$ws.Range("H113").Value = 1713.409
$ws.Range("I113").Value = 1599.7368
$ws.Range("J113").Value = 2433.3333
$ws.Range("K113").Value = 1599.7368
$ws.Range("L113").Value = 2433.3333
$ws.Range("M113").Value = 570.2632000000001
$ws.Range("N113").Value = -6773.3333
# Row 122
$ws.Range("H122").Value = 151813.86
$ws.Range("I122").Value = 204759.6
$ws.Range("J122").Value = 19449.5
$ws.Range("K122").Value = 614278.8
$ws.Range("L122").Value = 58348.5
$ws.Range("M122").Value = -611828.8
$ws.Range("N122").Value = -63248.5
# Row 132
$ws.Range("H132").Value = 2449.077
$ws.Range("I132").Value = 2449.0557
$ws.Range("K132").Value = 7347.1671
$ws.Range("M132").Value = -4817.1671

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 33
$ws.Range("H33").Value = 710.375
$ws.Range("I33").Value = 61
$ws.Range("J33").Value = 1100
$ws.Range("K33").Value = 366
$ws.Range("L33").Value = 6600
$ws.Range("M33").Value = -83
$ws.Range("N33").Value = -7166
# Row 59
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
# Row 68
$ws.Range("H68").Value = 1657.7142
$ws.Range("I68").Value = 1849.75
$ws.Range("J68").Value = 1401.6666
$ws.Range("K68").Value = 5549.25
$ws.Range("L68").Value = 4204.9998
$ws.Range("M68").Value = -4738.25
$ws.Range("N68").Value = -5826.9998
# Row 71
$ws.Range("H71").Value = 1657.7142
$ws.Range("I71").Value = 1849.75
$ws.Range("J71").Value = 1401.6666
$ws.Range("K71").Value = 16647.75
$ws.Range("L71").Value = 12614.9994
$ws.Range("M71").Value = -12591.75
$ws.Range("N71").Value = -20726.9994
# Row 86
$ws.Range("H86").Value = 462
$ws.Range("I86").Value = 462
$ws.Range("K86").Value = 1386
$ws.Range("M86").Value = -200
# Row 89
$ws.Range("H89").Value = 462
$ws.Range("I89").Value = 462
$ws.Range("K89").Value = 4158
$ws.Range("M89").Value = 1770

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 2575
$ws.Range("I113").Value = 2761.2
$ws.Range("J113").Value = 2388.8
$ws.Range("K113").Value = 2761.2
$ws.Range("L113").Value = 2388.8
$ws.Range("M113").Value = -591.1999999999998
$ws.Range("N113").Value = -6728.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 734.25
$ws.Range("I16").Value = 717.7143
$ws.Range("J16").Value = 850
$ws.Range("K16").Value = 717.7143
$ws.Range("L16").Value = 850
$ws.Range("M16").Value = -547.7143
$ws.Range("N16").Value = -1190
# Row 20
$ws.Range("H20").Value = 597975.5
$ws.Range("I20").Value = 9472.6
$ws.Range("K20").Value = 9472.6
$ws.Range("M20").Value = -9246.6
# Row 61
$ws.Range("H61").Value = 91885.80499999999
$ws.Range("I61").Value = 92709.664
$ws.Range("J61").Value = 88425.60000000001
$ws.Range("K61").Value = 92709.664
$ws.Range("L61").Value = 88425.60000000001
$ws.Range("M61").Value = -92507.664
$ws.Range("N61").Value = -88829.60000000001
# Row 82
$ws.Range("H82").Value = 5475.3335
$ws.Range("I82").Value = 5475.3335
$ws.Range("K82").Value = 5475.3335
$ws.Range("M82").Value = -5114.3335
# Row 85
$ws.Range("H85").Value = 5475.3335
$ws.Range("I85").Value = 5475.3335
$ws.Range("K85").Value = 5475.3335
$ws.Range("M85").Value = -4227.3335
# Row 93
$ws.Range("H93").Value = 29221
$ws.Range("I93").Value = 3948.7273
$ws.Range("K93").Value = 3948.7273
$ws.Range("M93").Value = -2700.7273
# Row 113
$ws.Range("H113").Value = 91885.80499999999
$ws.Range("I113").Value = 92709.664
$ws.Range("J113").Value = 88425.60000000001
$ws.Range("K113").Value = 92709.664
$ws.Range("L113").Value = 88425.60000000001
$ws.Range("M113").Value = -90539.664
$ws.Range("N113").Value = -92765.60000000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 7
$ws.Range("H7").Value = 9790
$ws.Range("I7").Value = 475
$ws.Range("K7").Value = 475
$ws.Range("M7").Value = -362
# Row 9
$ws.Range("H9").Value = 1500
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()
# Row 113
$ws.Range("H113").Value = 853.8333
$ws.Range("I113").Value = 780.75
$ws.Range("K113").Value = 2342.25
$ws.Range("M113").Value = -172.25

